$d = $word.ActiveDocument

function Find-ParagraphIndexLike($doc, $pattern) {
    $cnt = $doc.Paragraphs.Count
    for ($i = 1; $i -le $cnt; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text.Trim()
        if ($t -like $pattern) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------
# Hunk 1: "return e " paragraph -> "print e" (text "return " becomes
# "print ", and the trailing lone-space run after "e" is removed).
# ---------------------------------------------------------------
$idx1 = Find-ParagraphIndexLike $d "return e"
$p1 = $d.Paragraphs.Item($idx1)
$r1 = $p1.Range
$found1 = $r1.Find.Execute("return ", $true, $false, $false, $false, $false, $true, 1, $false, "print ", 2)

# Re-fetch the paragraph (text changed) and trim the trailing " " run
# that precedes the paragraph mark (i.e. delete the very last character
# before the end-of-paragraph mark, which is a lone space run).
$p1 = $d.Paragraphs.Item($idx1)
$r1 = $p1.Range
$e1 = $r1.End
# $e1 - 1 is the paragraph mark; the character right before that is the
# trailing space run to be removed entirely.
$rngTrailingSpace = $d.Range($e1 - 2, $e1 - 1)
$rngTrailingSpace.Delete()

# ---------------------------------------------------------------
# Hunk 2: "return total" paragraph -> split into an empty paragraph
# followed by a new paragraph starting with "print" (rest of the
# paragraph - " total" - is untouched).
# ---------------------------------------------------------------
$idx2 = Find-ParagraphIndexLike $d "return total"
$p2 = $d.Paragraphs.Item($idx2)
$r2 = $p2.Range
$s2 = $r2.Start
$rngReturn = $d.Range($s2, $s2 + 6)
$found2 = $rngReturn.Find.Execute("return", $true, $false, $false, $false, $false, $true, 1, $false, "^pprint", 2)

# ---------------------------------------------------------------
# Hunk 3: remove the blank paragraph and the page-break paragraph that
# sit between "...anywhere in between). " and "Functions:". Delete them
# one at a time (re-resolving the Range each time) because deleting a
# single Range spanning multiple paragraph marks only consumes one of
# them (Word leaves the last paragraph mark of the selection intact).
# ---------------------------------------------------------------
$idx3 = Find-ParagraphIndexLike $d "*anywhere in between*"

$pNext = $d.Paragraphs.Item($idx3 + 1)
$nStart = $pNext.Range.Start
$nEnd = $pNext.Range.End
$rngNext = $d.Range($nStart, $nEnd)
$rngNext.Delete()

$pNext = $d.Paragraphs.Item($idx3 + 1)
$nStart = $pNext.Range.Start
$nEnd = $pNext.Range.End
$rngNext = $d.Range($nStart, $nEnd)
$rngNext.Delete()

Write-Output "done"
